$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$c = $tcs.Colors(1)
$c.RGB = 255
Write-Host "after: $($tcs.Colors(1).RGB)"
